$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D hold price text that LOOKS numeric (e.g. "603.13"); Excel
# auto-converts such strings to real numbers on assignment. The source
# data must stay plain text (as it was authored), so each D-column write
# forces a Text number format before the assignment and then restores the
# cell to the (unstyled) "Normal" style so no stray formatting is left
# behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '63.675.96'
$ws.Range("E2").Value = '  +1.21%  '

Set-TextValue $ws.Range("D3") '3.282.87'
$ws.Range("E3").Value = '  +4.98%  '

$ws.Range("E4").Value = '  +0.16%  '

Set-TextValue $ws.Range("D5") '603.13'
$ws.Range("E5").Value = '  +2.30%  '

Set-TextValue $ws.Range("D6") '141.65'
$ws.Range("E6").Value = '  +4.31%  '

$ws.Range("E7").Value = '  +0.10%  '

Set-TextValue $ws.Range("D8") '3.278.61'
$ws.Range("E8").Value = '  +5.05%  '

Set-TextValue $ws.Range("D9") '0.518'
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("E10").Value = '  +3.23%  '

Set-TextValue $ws.Range("D11") '5.41'
$ws.Range("E11").Value = '  +3.26%  '

Set-TextValue $ws.Range("D12") '0.469'
$ws.Range("E12").Value = '  +3.47%  '

$ws.Range("E13").Value = '  +1.64%  '

Set-TextValue $ws.Range("D14") '34.47'
$ws.Range("E14").Value = '  +1.67%  '

Set-TextValue $ws.Range("D15") '3.826.60'
$ws.Range("E15").Value = '  +5.26%  '

$ws.Range("E16").Value = '  +1.04%  '

Set-TextValue $ws.Range("D17") '3.288.08'
$ws.Range("E17").Value = '  +5.43%  '

Set-TextValue $ws.Range("D18") '63.752.39'
$ws.Range("E18").Value = '  +1.27%  '

Set-TextValue $ws.Range("D19") '6.82'
$ws.Range("E19").Value = '  +3.27%  '

Set-TextValue $ws.Range("D20") '479.27'
$ws.Range("E20").Value = '  +2.29%  '

Set-TextValue $ws.Range("D21") '14.05'
$ws.Range("E21").Value = '  +0.05%  '

Set-TextValue $ws.Range("D22") '0.727'
$ws.Range("E22").Value = '  +4.74%  '

Set-TextValue $ws.Range("D23") '7.99'
$ws.Range("E23").Value = '  +4.47%  '

Set-TextValue $ws.Range("D24") '13.45'
$ws.Range("E24").Value = '  +4.39%  '

Set-TextValue $ws.Range("D25") '84.16'
$ws.Range("E25").Value = '  -1.12%  '

Set-TextValue $ws.Range("D26") '0.998'
$ws.Range("E26").Value = '  -0.22%  '

$ws.Range("E27").Value = '  +2.71%  '

Set-TextValue $ws.Range("D28") '7.29'
$ws.Range("E28").Value = '  +7.49%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D30") '2.17'
$ws.Range("E30").Value = '  +4.03%  '

$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D31") '8.10'
$ws.Range("E31").Value = '  +3.29%  '

Set-TextValue $ws.Range("D32") '28.80'
$ws.Range("E32").Value = '  +8.62%  '

$ws.Range("E33").Value = '  -3.11%  '

$ws.Range("E34").Value = '  +0.28%  '

Set-TextValue $ws.Range("D35") '1.09'
$ws.Range("E35").Value = '  +2.92%  '

Set-TextValue $ws.Range("D36") '5.95'
$ws.Range("E36").Value = '  +3.97%  '

Set-TextValue $ws.Range("D37") '52.98'
$ws.Range("E37").Value = '  +2.01%  '

Set-TextValue $ws.Range("D38") '0.0₃0741'
$ws.Range("E38").Value = '  +9.65%  '

Set-TextValue $ws.Range("D39") '0.0398'
$ws.Range("E39").Value = '  +3.58%  '

Set-TextValue $ws.Range("D40") '424.52'
$ws.Range("E40").Value = '  +2.24%  '

Set-TextValue $ws.Range("D41") '3.044.86'
$ws.Range("E41").Value = '  +4.93%  '

Set-TextValue $ws.Range("D42") '8.33'
$ws.Range("E42").Value = '  +2.05%  '

Set-TextValue $ws.Range("D43") '2.72'
$ws.Range("E43").Value = '  +2.26%  '

$ws.Range("E44").Value = '  +0.94%  '

Set-TextValue $ws.Range("D45") '0.264'
$ws.Range("E45").Value = '  +2.80%  '

Set-TextValue $ws.Range("D46") '2.18'
$ws.Range("E46").Value = '  +4.28%  '

Set-TextValue $ws.Range("D47") '26.13'
$ws.Range("E47").Value = '  +3.41%  '

$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("E49").Value = '  +2.28%  '

Set-TextValue $ws.Range("D50") '124.72'
$ws.Range("E50").Value = '  +3.55%  '

Set-TextValue $ws.Range("D51") '2.28'
$ws.Range("E51").Value = '  +2.17%  '
